$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: fix six pre-existing similarity values (recalculation fix) ---
# (coordinates below are in the ORIGINAL, pre-insert layout)
$ws.Cells.Item(3, 9).Value = 0.04
$ws.Cells.Item(9, 3).Value = 0.04
$ws.Cells.Item(9, 10).Value = 0.01
$ws.Cells.Item(10, 9).Value = 0.01
$ws.Cells.Item(28, 35).Value = 0.07000000000000001
$ws.Cells.Item(35, 28).Value = 0.07000000000000001
$ws.Cells.Item(31, 35).Value = 0.21
$ws.Cells.Item(35, 31).Value = 0.21
$ws.Cells.Item(32, 35).Value = 0.04
$ws.Cells.Item(35, 32).Value = 0.04
$ws.Cells.Item(35, 36).Value = 0.09
$ws.Cells.Item(36, 35).Value = 0.09

# --- Step 2: insert a new column for user 62 (becomes column AJ) ---
$ws.Columns.Item(36).Insert()

# --- Step 3: insert a new row for user 62 (becomes row 36) ---
$ws.Rows.Item(36).Insert()

# --- Step 4: copy header style onto the new column header cell, then set its value ---
$ws.Cells.Item(1, 1).Copy()
$ws.Cells.Item(1, 36).PasteSpecial(-4122)
$ws.Cells.Item(1, 36).Value = 62

# --- Step 5: copy row-label style onto the new row label cell, then set its value ---
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(36, 1).PasteSpecial(-4122)
$ws.Cells.Item(36, 1).Value = 62

# --- Step 6: fill the new row 36 (user 62 vs everyone), columns B..AN ---
$ws.Cells.Item(36, 2).Value = 0
$ws.Cells.Item(36, 3).Value = 0
$ws.Cells.Item(36, 4).Value = 0.18
$ws.Cells.Item(36, 5).Value = 0
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(36, 7).Value = 0
$ws.Cells.Item(36, 8).Value = 0
$ws.Cells.Item(36, 9).Value = 0
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 0
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).Value = 0
$ws.Cells.Item(36, 14).Value = 0
$ws.Cells.Item(36, 15).Value = 0
$ws.Cells.Item(36, 16).Value = 0
$ws.Cells.Item(36, 17).Value = 0
$ws.Cells.Item(36, 18).Value = 0.02
$ws.Cells.Item(36, 19).Value = 0
$ws.Cells.Item(36, 20).Value = 0
$ws.Cells.Item(36, 21).Value = 0
$ws.Cells.Item(36, 22).Value = 0
$ws.Cells.Item(36, 23).Value = 0
$ws.Cells.Item(36, 24).Value = 0
$ws.Cells.Item(36, 25).Value = 0
$ws.Cells.Item(36, 26).Value = 0
$ws.Cells.Item(36, 27).Value = 0
$ws.Cells.Item(36, 28).Value = 0
$ws.Cells.Item(36, 29).Value = 0
$ws.Cells.Item(36, 30).Value = 0
$ws.Cells.Item(36, 31).Value = 0
$ws.Cells.Item(36, 32).Value = 0
$ws.Cells.Item(36, 33).Value = 0
$ws.Cells.Item(36, 34).Value = 0.87
$ws.Cells.Item(36, 35).Value = 0
$ws.Cells.Item(36, 36).Value = 1
$ws.Cells.Item(36, 37).Value = 0
$ws.Cells.Item(36, 38).Value = 0.38
$ws.Cells.Item(36, 39).Value = 0
$ws.Cells.Item(36, 40).Value = 0

# --- Step 7: fill the new column AJ (user 62 vs everyone), rows 2..40 excluding 36 ---
$ws.Cells.Item(2, 36).Value = 0
$ws.Cells.Item(3, 36).Value = 0
$ws.Cells.Item(4, 36).Value = 0.18
$ws.Cells.Item(5, 36).Value = 0
$ws.Cells.Item(6, 36).Value = 0
$ws.Cells.Item(7, 36).Value = 0
$ws.Cells.Item(8, 36).Value = 0
$ws.Cells.Item(9, 36).Value = 0
$ws.Cells.Item(10, 36).Value = 0
$ws.Cells.Item(11, 36).Value = 0
$ws.Cells.Item(12, 36).Value = 0
$ws.Cells.Item(13, 36).Value = 0
$ws.Cells.Item(14, 36).Value = 0
$ws.Cells.Item(15, 36).Value = 0
$ws.Cells.Item(16, 36).Value = 0
$ws.Cells.Item(17, 36).Value = 0
$ws.Cells.Item(18, 36).Value = 0.02
$ws.Cells.Item(19, 36).Value = 0
$ws.Cells.Item(20, 36).Value = 0
$ws.Cells.Item(21, 36).Value = 0
$ws.Cells.Item(22, 36).Value = 0
$ws.Cells.Item(23, 36).Value = 0
$ws.Cells.Item(24, 36).Value = 0
$ws.Cells.Item(25, 36).Value = 0
$ws.Cells.Item(26, 36).Value = 0
$ws.Cells.Item(27, 36).Value = 0
$ws.Cells.Item(28, 36).Value = 0
$ws.Cells.Item(29, 36).Value = 0
$ws.Cells.Item(30, 36).Value = 0
$ws.Cells.Item(31, 36).Value = 0
$ws.Cells.Item(32, 36).Value = 0
$ws.Cells.Item(33, 36).Value = 0
$ws.Cells.Item(34, 36).Value = 0.87
$ws.Cells.Item(35, 36).Value = 0
$ws.Cells.Item(37, 36).Value = 0
$ws.Cells.Item(38, 36).Value = 0.38
$ws.Cells.Item(39, 36).Value = 0
$ws.Cells.Item(40, 36).Value = 0

$ws.Application.CutCopyMode = $false
